$wb = $excel.ActiveWorkbook

# --- About sheet: update "last updated" date (C1) from 1/3/2024 to 3/28/2024 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- FPIEBP sheet: update "hard coal" row (row 3) balancing priorities ---
# production=1 (was 3), imports=3 (was 2), exports=2 (was 1)
$wsData = $wb.Worksheets.Item("FPIEBP")
$wsData.Range("B3").Value = 1
$wsData.Range("C3").Value = 3
$wsData.Range("D3").Value = 2

# --- FPIEBP sheet: move the active selection from F4 to E3 ---
$wsData.Range("E3").Select()
